$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "a1": re-case the header row and append 3 data rows.
# ---------------------------------------------------------------------------
$wsA1 = $wb.Worksheets.Item("a1")

$wsA1.Range("A1").Value = "Id"
$wsA1.Range("B1").Value = "Serial"
$wsA1.Range("C1").Value = "Model"
$wsA1.Range("D1").Value = "Version"
$wsA1.Range("E1").Value = "Created"

# Serial (B) and Version (D) columns hold text, not numbers, even though the
# literal characters look numeric - format the cells as Text first so the
# values round-trip as strings instead of being auto-coerced to numbers.
$wsA1.Range("B2:B4").NumberFormat = "@"
$wsA1.Range("D2:D4").NumberFormat = "@"

$wsA1.Range("A2").Value = 1
$wsA1.Range("B2").Value = "123"
$wsA1.Range("C2").Value = "a1"
$wsA1.Range("D2").Value = "1"
$wsA1.Range("E2").Value = 0

$wsA1.Range("A3").Value = 3
$wsA1.Range("B3").Value = "113"
$wsA1.Range("C3").Value = "a1"
$wsA1.Range("D3").Value = "2"
$wsA1.Range("E3").Value = 0

$wsA1.Range("A4").Value = 4
$wsA1.Range("B4").Value = "113"
$wsA1.Range("C4").Value = "a1"
$wsA1.Range("D4").Value = "2"
$wsA1.Range("E4").Value = 0

# ---------------------------------------------------------------------------
# Sheet "A2": re-case the header row and append 1 data row.
# ---------------------------------------------------------------------------
$wsA2 = $wb.Worksheets.Item("A2")

$wsA2.Range("A1").Value = "Id"
$wsA2.Range("B1").Value = "Serial"
$wsA2.Range("C1").Value = "Model"
$wsA2.Range("D1").Value = "Version"
$wsA2.Range("E1").Value = "Created"

$wsA2.Range("B2").NumberFormat = "@"
$wsA2.Range("D2").NumberFormat = "@"

$wsA2.Range("A2").Value = 2
$wsA2.Range("B2").Value = "123"
$wsA2.Range("C2").Value = "A2"
$wsA2.Range("D2").Value = "1"
$wsA2.Range("E2").Value = 0
